# Update countries & provincias Spain
# Applies the data refresh captured in the commit "Update countries & provincias Spain":
#  - bumps the "datos actualizados" timestamp from 03:04 to 04:04
#  - refreshes numeric counters for several countries (Mexico, Corea del Sur, Australia,
#    Nicaragua/San Cristobal row, ...)
#  - Guatemala's case count overtook the row above it, so it (and the three rows that
#    used to sit between it and "Republica de Chipre") shift up/down by one rank
#  - Belice and Nueva Caledonia swap rank order as well

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (row 1) -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 10 de Mayo de 2020 a las 04:04"

# --- Mexico (row 21) -----------------------------------------------------------
$ws.Range("D21").Value = 21824
$ws.Range("E21").Value = 8283

# --- Corea del Sur (row 41) ----------------------------------------------------
$ws.Range("B41").Value = 10874
$ws.Range("C41").Value = 34
$ws.Range("D41").Value = 9610
$ws.Range("E41").Value = 1008

# --- Australia (row 53) --------------------------------------------------------
$ws.Range("B53").Value = 6931
$ws.Range("C53").Value = 2
$ws.Range("E53").Value = 699
$ws.Range("F53").Value = 17

# --- Guatemala climbs past Letonia / Kirguistan / Consejo Danes (rows 97-100) -
# Row 97 now holds Guatemala's (new, larger) figures ...
$ws.Range("A97").Value = "Guatemala"
$ws.Range("B97").Value = 967
$ws.Range("C97").Value = 67
$ws.Range("D97").Value = 104
$ws.Range("E97").Value = 839
$ws.Range("F97").Value = 5
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 24

# ... and the three rows that used to be above it each drop one rank, keeping
# their own previous figures.
$ws.Range("A98").Value = "Consejo Danes para los Refugiados"
$ws.Range("B98").Value = 937
$ws.Range("C98").Value = 0
$ws.Range("D98").Value = 130
$ws.Range("E98").Value = 768
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 39

$ws.Range("A99").Value = "Kirguistan"
$ws.Range("B99").Value = 931
$ws.Range("C99").Value = 0
$ws.Range("D99").Value = 658
$ws.Range("E99").Value = 261
$ws.Range("F99").Value = 13
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 12

$ws.Range("A100").Value = "Letonia"
$ws.Range("B100").Value = 930
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 464
$ws.Range("E100").Value = 448
$ws.Range("F100").Value = 2
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 18

# --- Belice / Nueva Caledonia swap rank (rows 192-193) ------------------------
$ws.Range("A192").Value = "Nueva Caledonia"
$ws.Range("B192").Value = 18
$ws.Range("C192").Value = 0
$ws.Range("D192").Value = 18
$ws.Range("E192").Value = 0
$ws.Range("F192").Value = 0
$ws.Range("G192").Value = 0
$ws.Range("H192").Value = 0

$ws.Range("A193").Value = "Belice"
$ws.Range("B193").Value = 18
$ws.Range("C193").Value = 0
$ws.Range("D193").Value = 16
$ws.Range("E193").Value = 0
$ws.Range("F193").Value = 0
$ws.Range("G193").Value = 0
$ws.Range("H193").Value = 2

# --- San Cristobal y Nieves (row 201) -----------------------------------------
$ws.Range("D201").Value = 14
$ws.Range("E201").Value = 1
